$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 87: B87 should become a real number (3) instead of text "3".
# Clear and re-set as numeric value.
$ws.Cells.Item(87, 2).Value = 3

# Row 88: new row of data appended after row 87.
$ws.Cells.Item(88, 1).Value = "Ruilin"
$ws.Cells.Item(88, 2).NumberFormat = "@"
$ws.Cells.Item(88, 2).Value = "3"
$ws.Cells.Item(88, 2).Style = "Normal"
$ws.Cells.Item(88, 3).Value = "无"
$ws.Cells.Item(88, 4).Value = "DFT"
$ws.Cells.Item(88, 5).Value = "WRI"
$ws.Cells.Item(88, 6).Value = "4d925051-98cd-4c47-ad57-530c84f23548"
$ws.Cells.Item(88, 7).Value = "B1QRgziT-_annotated.xlsx"
$ws.Cells.Item(88, 8).Value = "I don't think this paper explains the importance of its results nearly enough and I'm concerned that it may not be obvious what a breakthrough it is just from skimming the abstract."
